$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 8 new rows at row 3 (pushes old rows 3-4 down to 11-12)
$ws.Rows("3:10").Insert()

# Populate the newly inserted rows with the related-companies data
# Row 3
$ws.Cells.Item(3, 1).Value = 'Relacionado (Cascante)'
$ws.Cells.Item(3, 2).Value = 'Inyecciones Plásticas Mecacontrol'
$ws.Cells.Item(3, 3).Value = 'Carr. Tudela - Tarazona, 10, 31520 Cascante, Navarra, España'
$ws.Cells.Item(3, 4).Value = '+34 948 84 45 21'
$ws.Cells.Item(3, 5).Value = 'http://www.inyeccionesmecacontrol.com/'
$ws.Cells.Item(3, 6).Value = 'https://www.google.com/maps/place/Inyecciones+Pl%C3%A1sticas+Mecacontrol/@41.9894472,-1.6827796,17z/data=!4m10!1m2!2m1!1sEmpresas+en+Cascante+Espa%C3%B1a!3m6!1s0xd5a51fdf53dea27:0x9823af5de2de00a8!8m2!3d41.9894472!4d-1.6782735!15sChxFbXByZXNhcyBlbiBDYXNjYW50ZSBFc3Bhw7FhkgEhcGxhc3RpY19pbmplY3Rpb25fbW9sZGluZ19zZXJ2aWNlqgFUEAEqDCIIZW1wcmVzYXMoDjIgEAEiHPN-cHidLlggMybfITUEhVsbNWgSxpnAvZ-OxigyIBACIhxlbXByZXNhcyBlbiBjYXNjYW50ZSBlc3Bhw7Fh4AEA!16s%2Fg%2F1q5bnr7lf?entry=ttu&g_ep=EgoyMDI1MDYwNC4wIKXMDSoASAFQAw%3D%3D'
$ws.Cells.Item(3, 7).Value = 'Empresa de moldeo por inyección de plástico'

# Row 4
$ws.Cells.Item(4, 1).Value = 'Relacionado (Cascante)'
$ws.Cells.Item(4, 2).Value = 'Vibracoustic Cascante'
$ws.Cells.Item(4, 3).Value = 'Ctra. Tudela-Tarazona, s/n, 31520 Cascante, Navarra, España'
$ws.Cells.Item(4, 4).Value = '+34 948 84 45 06'
$ws.Cells.Item(4, 5).Value = 'https://www.vibracoustic.com/'
$ws.Cells.Item(4, 6).Value = 'https://www.google.com/maps/place/Vibracoustic+Cascante/@41.9894472,-1.6827796,17z/data=!4m10!1m2!2m1!1sEmpresas+en+Cascante+Espa%C3%B1a!3m6!1s0xd5a51fe3f7664cb:0xf61aec36796f3de!8m2!3d41.990382!4d-1.6777962!15sChxFbXByZXNhcyBlbiBDYXNjYW50ZSBFc3Bhw7FhkgEXYXV0b19wYXJ0c19tYW51ZmFjdHVyZXKqAVQQASoMIghlbXByZXNhcygOMiAQASIc835weJ0uWCAzJt8hNQSFWxs1aBLGmcC9n47GKDIgEAIiHGVtcHJlc2FzIGVuIGNhc2NhbnRlIGVzcGHDsWHgAQA!16s%2Fg%2F1trszbgr?entry=ttu&g_ep=EgoyMDI1MDYwNC4wIKXMDSoASAFQAw%3D%3D'
$ws.Cells.Item(4, 7).Value = 'Fabricante de repuestos para automóviles'

# Row 5
$ws.Cells.Item(5, 1).Value = 'Relacionado (Cascante)'
$ws.Cells.Item(5, 2).Value = 'Industrias Areso'
$ws.Cells.Item(5, 3).Value = 'Pol. Industrial El Parral, Nave 8-9, 31520 Cascante, Navarra, España'
$ws.Cells.Item(5, 4).Value = '+34 948 85 15 52'
$ws.Cells.Item(5, 5).Value = 'http://www.aresoguantes.com/'
$ws.Cells.Item(5, 6).Value = 'https://www.google.com/maps/place/Industrias+Areso/@41.9894472,-1.6827796,17z/data=!4m10!1m2!2m1!1sEmpresas+en+Cascante+Espa%C3%B1a!3m6!1s0xd5a51fdf53dea27:0x45792b426f5264cd!8m2!3d41.9917081!4d-1.6776739!15sChxFbXByZXNhcyBlbiBDYXNjYW50ZSBFc3Bhw7FhWh4iHGVtcHJlc2FzIGVuIGNhc2NhbnRlIGVzcGHDsWGSARxwcm90ZWN0aXZlX2Nsb3RoaW5nX3N1cHBsaWVymgEkQ2hkRFNVaE5NRzluUzBWSlEwRm5TVU4yY0hGRFdEUkJSUkFCqgFUEAEqDCIIZW1wcmVzYXMoDjIgEAEiHPN-cHidLlggMybfITUEhVsbNWgSxpnAvZ-OxigyIBACIhxlbXByZXNhcyBlbiBjYXNjYW50ZSBlc3Bhw7Fh4AEA-gEECAAQOw!16s%2Fg%2F12qghm7tl?entry=ttu&g_ep=EgoyMDI1MDYwNC4wIKXMDSoASAFQAw%3D%3D'

# Row 6
$ws.Cells.Item(6, 1).Value = 'Relacionado (Cascante)'
$ws.Cells.Item(6, 2).Value = 'Edmar S A'
$ws.Cells.Item(6, 3).Value = 'Carr. Tudela - Tarazona, 23, 31520 Cascante, Navarra, España'
$ws.Cells.Item(6, 4).Value = '+34 948 85 02 83'
$ws.Cells.Item(6, 6).Value = 'https://www.google.com/maps/place/Edmar+S+A/@41.9962921,-1.6788175,17z/data=!4m10!1m2!2m1!1sEmpresas+en+Cascante+Espa%C3%B1a!3m6!1s0xd5a4e02ce7d0713:0xa533a4393fce17f3!8m2!3d41.9962921!4d-1.6743114!15sChxFbXByZXNhcyBlbiBDYXNjYW50ZSBFc3Bhw7FhkgEdaW5kdXN0cmlhbF9lcXVpcG1lbnRfc3VwcGxpZXKqAVQQASoMIghlbXByZXNhcygOMiAQASIc835weJ0uWCAzJt8hNQSFWxs1aBLGmcC9n47GKDIgEAIiHGVtcHJlc2FzIGVuIGNhc2NhbnRlIGVzcGHDsWHgAQA!16s%2Fg%2F1z44b4sz9?entry=ttu&g_ep=EgoyMDI1MDYwNC4wIKXMDSoASAFQAw%3D%3D'
$ws.Cells.Item(6, 7).Value = 'Empresa de suministros industriales'

# Row 7
$ws.Cells.Item(7, 1).Value = 'Relacionado (Cascante)'
$ws.Cells.Item(7, 2).Value = 'Galipienzo'
$ws.Cells.Item(7, 3).Value = 'C. Vía Romana, 0, 31520 Cascante, Navarra, España'
$ws.Cells.Item(7, 4).Value = '+34 948 85 16 66'
$ws.Cells.Item(7, 5).Value = 'https://galipienzo.es/'
$ws.Cells.Item(7, 6).Value = 'https://www.google.com/maps/place/Galipienzo/@42.0107959,-1.6882863,17z/data=!4m10!1m2!2m1!1sEmpresas+en+Cascante+Espa%C3%B1a!3m6!1s0xd5a51e23bf46c7f:0xfb0b0ad90ea41d5e!8m2!3d42.0107959!4d-1.6837802!15sChxFbXByZXNhcyBlbiBDYXNjYW50ZSBFc3Bhw7FhkgEOc2xhdWdodGVyaG91c2WqAVQQASoMIghlbXByZXNhcygOMiAQASIc835weJ0uWCAzJt8hNQSFWxs1aBLGmcC9n47GKDIgEAIiHGVtcHJlc2FzIGVuIGNhc2NhbnRlIGVzcGHDsWHgAQA!16s%2Fg%2F12qggw6dh?entry=ttu&g_ep=EgoyMDI1MDYwNC4wIKXMDSoASAFQAw%3D%3D'
$ws.Cells.Item(7, 7).Value = 'Matadero'

# Row 8
$ws.Cells.Item(8, 1).Value = 'Relacionado (Cascante)'
$ws.Cells.Item(8, 2).Value = 'Bodegas Malón de Echaide'
$ws.Cells.Item(8, 3).Value = 'Carr. Tudela - Tarazona, 33, 31520 Cascante, Navarra, España'
$ws.Cells.Item(8, 4).Value = '+34 948 85 14 11'
$ws.Cells.Item(8, 5).Value = 'http://www.malondeechaide.com/'
$ws.Cells.Item(8, 6).Value = 'https://www.google.com/maps/place/Bodegas+Mal%C3%B3n+de+Echaide/@41.9937281,-1.679878,17z/data=!4m10!1m2!2m1!1sEmpresas+en+Cascante+Espa%C3%B1a!3m6!1s0xd5a4e029a58e90f:0x6b8273c2b56a20b5!8m2!3d41.9937281!4d-1.6753719!15sChxFbXByZXNhcyBlbiBDYXNjYW50ZSBFc3Bhw7FhWh4iHGVtcHJlc2FzIGVuIGNhc2NhbnRlIGVzcGHDsWGSAQZ3aW5lcnmaASRDaGREU1VoTk1HOW5TMFZKUTBGblNVUmlhR05FYkdsUlJSQUKqAVQQASoMIghlbXByZXNhcygOMiAQASIc835weJ0uWCAzJt8hNQSFWxs1aBLGmcC9n47GKDIgEAIiHGVtcHJlc2FzIGVuIGNhc2NhbnRlIGVzcGHDsWHgAQD6AQUIlgIQEg!16s%2Fg%2F1tdc6kj8?entry=ttu&g_ep=EgoyMDI1MDYwNC4wIKXMDSoASAFQAw%3D%3D'
$ws.Cells.Item(8, 7).Value = 'Bodega'

# Row 9
$ws.Cells.Item(9, 1).Value = 'Relacionado (Cascante)'
$ws.Cells.Item(9, 2).Value = 'Centro Termolúdico Cascante'
$ws.Cells.Item(9, 3).Value = 'C. Fundacion Fuentes Dutor, s/n, 31520 Cascante, Navarra, España'
$ws.Cells.Item(9, 4).Value = '+34 948 84 45 38'
$ws.Cells.Item(9, 5).Value = 'http://www.termoludicocascante.es/'
$ws.Cells.Item(9, 6).Value = 'https://www.google.com/maps/place/Centro+Termol%C3%BAdico+Cascante/@41.9931,-1.6918171,17z/data=!4m10!1m2!2m1!1sEmpresas+en+Cascante+Espa%C3%B1a!3m6!1s0xd5a51f06019eabf:0xa4ee85a48f3119fc!8m2!3d41.9931!4d-1.687311!15sChxFbXByZXNhcyBlbiBDYXNjYW50ZSBFc3Bhw7FhWh4iHGVtcHJlc2FzIGVuIGNhc2NhbnRlIGVzcGHDsWGSARNzcGFfYW5kX2hlYWx0aF9jbHVimgEkQ2hkRFNVaE5NRzluUzBWSlEwRm5TVVIxYVdWbWJ6ZEJSUkFCqgFUEAEqDCIIZW1wcmVzYXMoDjIgEAEiHPN-cHidLlggMybfITUEhVsbNWgSxpnAvZ-OxigyIBACIhxlbXByZXNhcyBlbiBjYXNjYW50ZSBlc3Bhw7Fh4AEA-gEECAsQRg!16s%2Fg%2F1v_z5nbn?entry=ttu&g_ep=EgoyMDI1MDYwNC4wIKXMDSoASAFQAw%3D%3D'
$ws.Cells.Item(9, 7).Value = 'Spa y gimnasio'

# Row 10
$ws.Cells.Item(10, 1).Value = 'Relacionado (Cascante)'
$ws.Cells.Item(10, 2).Value = 'EADEC · Etiquetas Adhesivas Ecológicas'
$ws.Cells.Item(10, 3).Value = 'Polígono el Parral, 5, 31520 Cascante, Navarra, España'
$ws.Cells.Item(10, 4).Value = '+34 948 85 08 35'
$ws.Cells.Item(10, 5).Value = 'https://eadec.es/'
$ws.Cells.Item(10, 6).Value = 'https://www.google.com/maps/place/EADEC+%C2%B7+Etiquetas+Adhesivas+Ecol%C3%B3gicas/@41.9908239,-1.6824224,17z/data=!4m10!1m2!2m1!1sEmpresas+en+Cascante+Espa%C3%B1a!3m6!1s0xd5a4e027ecb5463:0xfcb68d187c69f39c!8m2!3d41.9908239!4d-1.6779163!15sChxFbXByZXNhcyBlbiBDYXNjYW50ZSBFc3Bhw7FhWh4iHGVtcHJlc2FzIGVuIGNhc2NhbnRlIGVzcGHDsWGSARRzdGlja2VyX21hbnVmYWN0dXJlcqoBVBABKgwiCGVtcHJlc2FzKA4yIBABIhzzfnB4nS5YIDMm3yE1BIVbGzVoEsaZwL2fjsYoMiAQAiIcZW1wcmVzYXMgZW4gY2FzY2FudGUgZXNwYcOxYeABAA!16s%2Fg%2F1tgzcp6s?entry=ttu&g_ep=EgoyMDI1MDYwNC4wIKXMDSoASAFQAw%3D%3D'
$ws.Cells.Item(10, 7).Value = 'Impresora de etiquetas personalizadas'

# Apply the light-blue highlight fill to the new rows (A3:G10)
$ws.Range("A3:G10").Interior.Color = 16770508

Write-Host "done"